$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list_study4_try_counter1")
$ws.Range("H2:H115").Value = "NaN"
$null = $ws.Range("H2:H115").Select()
